# "Generate Report for Handoff"
# Updates the localization-status report: status text, priority, timestamps,
# an error-detail message for a stale handback file, and a couple of
# column-width tweaks (to make room for the new long text in col P / E-F).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7328e9d3808a083ca292e1abb5cfcc29efc2f11f/e2e/a36f76a9-7d0b-4aec-9452-d115bc0d5281.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/147936ad1aa1a85ca48d8c0c46849d78def76fdf/e2e/a36f76a9-7d0b-4aec-9452-d115bc0d5281.md."

# ---------------------------------------------------------------------------
# Overview sheet: status moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and the generate-date timestamp advances.
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-11-15 18:04:20"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-11-15 18:04:20"

# ---------------------------------------------------------------------------
# zh-cn sheet: Status + Priority + Latest Handoff Datetime refresh, and the
# second file now reports a stale-handback Error Detail.
# ---------------------------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-11-15 18:04:07"

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-11-15 18:04:07"
$wsZhCn.Range("P3").Value = $errorDetail

# ---------------------------------------------------------------------------
# de-de sheet: same Status + Priority refresh, same Error Detail on row 3.
# Its Latest Handoff Datetime shared the same underlying string as the
# Overview sheet's generate-date ("2016-11-15 18:02:08"), so it advances to
# the same new timestamp too.
# ---------------------------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-11-15 18:04:20"

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-11-15 18:04:20"
$wsDeDe.Range("P3").Value = $errorDetail

# ---------------------------------------------------------------------------
# Column width tweaks. Excel snaps ColumnWidth to whole-pixel boundaries, so
# feed it the character width whose pixel rounding lands closest to the
# target (17.2159881591797 -> 17.1666..., 13.7470531463623 -> 40 exactly).
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
